$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "x" mark from E12 to D12
$ws.Range("E12").Value = ""
$ws.Range("D12").Value = "x"

# Add an "x" mark in D13
$ws.Range("D13").Value = "x"

# Update the view: scroll so A4 is the top-left visible cell, and set the active selection to D13
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D13").Select()
